$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.871.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.522.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.520.02'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  -2.50%  '
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.968.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.891.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.525.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.72%  '
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("E37").Value = '  -3.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.65'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.816'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("E41").Value = '  -1.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '285.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '132.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.81%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.997'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0929'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.93%  '
